$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 330.82352
$ws.Range("I42").Value = 48.444443
$ws.Range("J42").Value = 648.5
$ws.Range("K42").Value = 145.333329
$ws.Range("L42").Value = 1945.5
$ws.Range("M42").Value = 84.66667100000001
$ws.Range("N42").Value = -2405.5
$ws.Range("H58").Value = 99.71429000000001
$ws.Range("I58").Value = 99.71429000000001
$ws.Range("K58").Value = 299.14287
$ws.Range("M58").Value = -149.14287
$ws.Range("H113").Value = 3599.7
$ws.Range("J113").Value = 3833.1667
$ws.Range("L113").Value = 3833.1667
$ws.Range("N113").Value = -10341.1667
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
$ws.Range("H131").Value = 5141.7144
$ws.Range("I131").Value = 4332
$ws.Range("J131").Value = 10000
$ws.Range("K131").Value = 12996
$ws.Range("L131").Value = 30000
$ws.Range("M131").Value = -7956
$ws.Range("N131").Value = -40080
$ws.Range("H132").Value = 8136.51
$ws.Range("I132").Value = 2488.0222
$ws.Range("K132").Value = 7464.0666
$ws.Range("M132").Value = -4934.0666
$ws.Range("H137").Value = 3379.6667
$ws.Range("I137").Value = 3858.389
$ws.Range("K137").Value = 11575.167
$ws.Range("M137").Value = -9025.167000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6342.029
$ws.Range("I61").Value = 6467.875
$ws.Range("K61").Value = 6467.875
$ws.Range("M61").Value = -6255.875
$ws.Range("H63").Value = 2713.5715
$ws.Range("I63").Value = 2499
$ws.Range("J63").Value = 3250
$ws.Range("K63").Value = 2499
$ws.Range("L63").Value = 3250
$ws.Range("N63").Value = -4622
$ws.Range("M63").Value = -1813
$ws.Range("H66").Value = 2713.5715
$ws.Range("I66").Value = 2499
$ws.Range("J66").Value = 3250
$ws.Range("K66").Value = 12495
$ws.Range("L66").Value = 16250
$ws.Range("N66").Value = -23114
$ws.Range("M66").Value = -9063
$ws.Range("H74").Value = 1973.4572
$ws.Range("I74").Value = 754.1070999999999
$ws.Range("K74").Value = 754.1070999999999
$ws.Range("M74").Value = 119.8929000000001
$ws.Range("H77").Value = 1973.4572
$ws.Range("I77").Value = 754.1070999999999
$ws.Range("K77").Value = 3770.5355
$ws.Range("M77").Value = 597.4645
$ws.Range("H132").Value = 3109.7358
$ws.Range("I132").Value = 3016.0435
$ws.Range("K132").Value = 9048.130500000001
$ws.Range("M132").Value = -6518.130500000001
$ws.Range("H136").Value = 6342.029
$ws.Range("I136").Value = 6467.875
$ws.Range("K136").Value = 19403.625
$ws.Range("M136").Value = -16853.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H133").Value = 99776.664
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 99776.664
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 99776.664
$ws.Range("M133").ClearContents()
$ws.Range("N133").Value = -109896.664
$ws.Range("H134").Value = 6075.3657
$ws.Range("I134").Value = 3063.2632
$ws.Range("J134").Value = 8676.727999999999
$ws.Range("K134").Value = 9189.7896
$ws.Range("L134").Value = 26030.184
$ws.Range("M134").Value = -6654.7896
$ws.Range("N134").Value = -31100.184

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 5378
$ws.Range("H126").Value = 5378

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 4344.5
$ws.Range("I118").Value = 355.66666
$ws.Range("K118").Value = 1066.99998
$ws.Range("M118").Value = 176.0000199999999
$ws.Range("H134").Value = 1718.8125
$ws.Range("I134").Value = 1718.8125
$ws.Range("K134").Value = 5156.4375
$ws.Range("M134").Value = -86.4375
$ws.Range("H138").Value = 6477627
$ws.Range("I138").Value = 1120518.1
$ws.Range("J138").Value = 12504375
$ws.Range("K138").Value = 3361554.3
$ws.Range("L138").Value = 37513125
$ws.Range("M138").Value = -3356414.3
$ws.Range("N138").Value = -37523405
$ws.Range("H139").Value = 5559361
$ws.Range("I139").Value = 2373.75
$ws.Range("J139").Value = 10004951
$ws.Range("K139").Value = 7121.25
$ws.Range("L139").Value = 30014853
$ws.Range("M139").Value = -1981.25
$ws.Range("N139").Value = -30025133
$ws.Range("H140").Value = 2338.5938
$ws.Range("I140").Value = 2030.9259
$ws.Range("J140").Value = 4000
$ws.Range("K140").Value = 6092.7777
$ws.Range("L140").Value = 12000
$ws.Range("M140").Value = -912.7776999999996
$ws.Range("N140").Value = -22360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 32166.666
$ws.Range("J123").Value = 32166.666
$ws.Range("L123").Value = 32166.666
$ws.Range("N123").Value = -37066.666
$ws.Range("H132").Value = 5558888
$ws.Range("I132").Value = 6669876.5
$ws.Range("J132").Value = 3946.5
$ws.Range("K132").Value = 20009629.5
$ws.Range("L132").Value = 11839.5
$ws.Range("M132").Value = -20007099.5
$ws.Range("N132").Value = -16899.5
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4044.4482
$ws.Range("I7").Value = 3115.4211
$ws.Range("J7").Value = 5809.6
$ws.Range("K7").Value = 3115.4211
$ws.Range("L7").Value = 5809.6
$ws.Range("M7").Value = -3003.4211
$ws.Range("N7").Value = -6033.6
$ws.Range("H126").Value = 4044.4482
$ws.Range("I126").Value = 3115.4211
$ws.Range("J126").Value = 5809.6
$ws.Range("K126").Value = 9346.263300000001
$ws.Range("L126").Value = 17428.8
$ws.Range("M126").Value = -6876.263300000001
$ws.Range("N126").Value = -22368.8
$ws.Range("H132").Value = 2794.0735
$ws.Range("I132").Value = 2732.9814
$ws.Range("K132").Value = 8198.9442
$ws.Range("M132").Value = -5668.9442
$ws.Range("H136").Value = 5096.846
$ws.Range("I136").Value = 4811.9375
$ws.Range("K136").Value = 14435.8125
$ws.Range("M136").Value = -11885.8125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 50000
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("H96").Value = 2881.5715
$ws.Range("I96").Value = 2099.875
$ws.Range("J96").Value = 3923.8333
$ws.Range("K96").Value = 2099.875
$ws.Range("L96").Value = 3923.8333
$ws.Range("M96").Value = -726.875
$ws.Range("N96").Value = -6669.8333
$ws.Range("H100").Value = 542.8
$ws.Range("I100").Value = 569.7778
$ws.Range("K100").Value = 1139.5556
$ws.Range("M100").Value = -598.5555999999999
$ws.Range("H109").Value = 89618.336
$ws.Range("J109").Value = 89618.336
$ws.Range("L109").Value = 89618.336
$ws.Range("N109").Value = -92392.336
$ws.Range("H132").Value = 2392.923
$ws.Range("I132").Value = 2463.4211
$ws.Range("K132").Value = 7390.263300000001
$ws.Range("M132").Value = -4860.263300000001
$ws.Range("H135").Value = 70166.664
$ws.Range("J135").Value = 70166.664
$ws.Range("L135").Value = 70166.664
$ws.Range("N135").Value = -80306.664
$ws.Range("H136").Value = 13114.223
$ws.Range("I136").Value = 17005.691
$ws.Range("J136").Value = 2996.4
$ws.Range("K136").Value = 51017.073
$ws.Range("L136").Value = 8989.200000000001
$ws.Range("M136").Value = -48467.073
$ws.Range("N136").Value = -14089.2
